$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-01-01 Wednesday" "2025-01-02 Thursday"

Replace-Text "67×65=" "16×59="
Replace-Text "76×60=" "79×56="
Replace-Text "84×70=" "68×60="
Replace-Text "18×44=" "71×35="
Replace-Text "71×43=" "17×66="
Replace-Text "80×91=" "42×88="
Replace-Text "60×97=" "60×50="
Replace-Text "82×82=" "98×67="
Replace-Text "79×82=" "79×61="
Replace-Text "87×43=" "68×77="
Replace-Text "23×32=" "97×41="
Replace-Text "83×87=" "70×95="
Replace-Text "69×58=" "82×75="
Replace-Text "83×76=" "23×35="
Replace-Text "11×20=" "35×77="
Replace-Text "21×14=" "94×25="
Replace-Text "11×52=" "44×30="
Replace-Text "31×22=" "79×18="
Replace-Text "51×46=" "81×88="
Replace-Text "57×56=" "32×92="
Replace-Text "64×29=" "54×12="
Replace-Text "41×52=" "53×66="
Replace-Text "97×15=" "31×29="
Replace-Text "47×98=" "74×89="
Replace-Text "93×11=" "73×92="
